$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the total_measure values for rows 17 and 18 (column D)
$ws.Range("D17").Value = 480
$ws.Range("D18").Value = 240

# Update the active selection to match the new edit location
$ws.Range("D18").Select()
